# Atualizacao de bases das ligas, do dia: 21-04-2024 as 14:32
#
# Four new match rows are added to the "Hungary NB I" results table. The
# sheet's id column (A) is a contiguous 0-based row index, so the three
# rows that used to be the last rows in the sheet (old rows 170-172, ids
# 168-170, not yet played at the time of the previous snapshot) are moved
# down to rows 174-176 (ids 172-174) unchanged, and rows 170-173 (ids
# 168-171) are filled with fresh data for matches that have since been
# played (including their final-score FTHG/FTAG/FTR columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the old rows 170:172 down to 174:176 ----------------
# Copy formats, then values, for just the populated column blocks (A:G and
# K:AA - H/I/J/AB/AC were not populated on these rows) so we don't use
# EntireRow.Insert (which would otherwise synthesize a new, unused style).
$ws.Range("A170:G172").Copy()
$ws.Range("A174:G176").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K170:AA172").Copy()
$ws.Range("K174:AA176").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A170:G172").Copy()
$ws.Range("A174:G176").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("K170:AA172").Copy()
$ws.Range("K174:AA176").PasteSpecial(-4163)  # xlPasteValues

# Renumber the relocated rows' id column (A) so the index stays contiguous.
$ws.Range("A174").Value = 172
$ws.Range("A175").Value = 173
$ws.Range("A176").Value = 174

# --- Step 2: give the brand-new row 173 the same formatting as row 172 ----
$ws.Range("A172:AC172").Copy()
$ws.Range("A173:AC173").PasteSpecial(-4122)  # xlPasteFormats

# --- Step 3: write the new match data into rows 170-173 -------------------

# Row 170: Puskas Academy vs Paksi (2024-04-20 09:30)
$ws.Range("A170").Value = 168
$ws.Range("B170").Value = 6818365
$ws.Range("C170").Value = "Hungary NB I"
$ws.Range("D170").Value = "Hungary NB I"
$ws.Range("E170").Value = 45402.39583333334
$ws.Range("F170").Value = "Puskas Academy"
$ws.Range("G170").Value = "Paksi"
$ws.Range("H170").Value = 5
$ws.Range("I170").Value = 0
$ws.Range("J170").Value = "H"
$ws.Range("K170").Value = 2
$ws.Range("L170").Value = 3.4
$ws.Range("M170").Value = 3.3
$ws.Range("N170").Value = 2.05
$ws.Range("O170").Value = 3.25
$ws.Range("P170").Value = 3.3
$ws.Range("Q170").Value = -0.25
$ws.Range("R170").Value = 1.85
$ws.Range("S170").Value = 2
$ws.Range("T170").Value = 2.75
$ws.Range("U170").Value = 1.825
$ws.Range("V170").Value = 2.025
$ws.Range("W170").Value = 1.05
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = 0.8500000000000001
$ws.Range("AA170").Value = -1
$ws.Range("AB170").Value = 0.825
$ws.Range("AC170").Value = -1

# Row 171: Ferencvarosi TC vs Kisvarda FC (2024-04-20 12:00)
$ws.Range("A171").Value = 169
$ws.Range("B171").Value = 6818362
$ws.Range("C171").Value = "Hungary NB I"
$ws.Range("D171").Value = "Hungary NB I"
$ws.Range("E171").Value = 45402.5
$ws.Range("F171").Value = "Ferencvarosi TC"
$ws.Range("G171").Value = "Kisvarda FC"
$ws.Range("H171").Value = 0
$ws.Range("I171").Value = 0
$ws.Range("J171").Value = "D"
$ws.Range("K171").Value = 1.25
$ws.Range("L171").Value = 5.5
$ws.Range("M171").Value = 9
$ws.Range("N171").Value = 1.125
$ws.Range("O171").Value = 7
$ws.Range("P171").Value = 17
$ws.Range("Q171").Value = -2.25
$ws.Range("R171").Value = 1.95
$ws.Range("S171").Value = 1.9
$ws.Range("T171").Value = 3.25
$ws.Range("U171").Value = 1.975
$ws.Range("V171").Value = 1.875
$ws.Range("W171").Value = -1
$ws.Range("X171").Value = 6
$ws.Range("Y171").Value = -1
$ws.Range("Z171").Value = -1
$ws.Range("AA171").Value = 0.8999999999999999
$ws.Range("AB171").Value = -1
$ws.Range("AC171").Value = 0.875

# Row 172: Diosgyori VTK vs Debreceni VSC (2024-04-20 14:30)
$ws.Range("A172").Value = 170
$ws.Range("B172").Value = 6818364
$ws.Range("C172").Value = "Hungary NB I"
$ws.Range("D172").Value = "Hungary NB I"
$ws.Range("E172").Value = 45402.60416666666
$ws.Range("F172").Value = "Diosgyori VTK"
$ws.Range("G172").Value = "Debreceni VSC"
$ws.Range("H172").Value = 5
$ws.Range("I172").Value = 3
$ws.Range("J172").Value = "H"
$ws.Range("K172").Value = 2.5
$ws.Range("L172").Value = 3.2
$ws.Range("M172").Value = 2.625
$ws.Range("N172").Value = 2.7
$ws.Range("O172").Value = 3.25
$ws.Range("P172").Value = 2.4
$ws.Range("Q172").Value = 0
$ws.Range("R172").Value = 2.05
$ws.Range("S172").Value = 1.8
$ws.Range("T172").Value = 2.75
$ws.Range("U172").Value = 1.975
$ws.Range("V172").Value = 1.875
$ws.Range("W172").Value = 1.7
$ws.Range("X172").Value = -1
$ws.Range("Y172").Value = -1
$ws.Range("Z172").Value = 1.05
$ws.Range("AA172").Value = -1
$ws.Range("AB172").Value = 0.9750000000000001
$ws.Range("AC172").Value = -1

# Row 173: Kecskemeti TE vs Mezokovesd Zsory (2024-04-21 09:05)
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 6818367
$ws.Range("C173").Value = "Hungary NB I"
$ws.Range("D173").Value = "Hungary NB I"
$ws.Range("E173").Value = 45403.37847222222
$ws.Range("F173").Value = "Kecskemeti TE"
$ws.Range("G173").Value = "Mezokovesd Zsory"
$ws.Range("H173").Value = 2
$ws.Range("I173").Value = 1
$ws.Range("J173").Value = "H"
$ws.Range("K173").Value = 1.727
$ws.Range("L173").Value = 3.5
$ws.Range("M173").Value = 4.333
$ws.Range("N173").Value = 1.4
$ws.Range("O173").Value = 4.2
$ws.Range("P173").Value = 7
$ws.Range("Q173").Value = -1.25
$ws.Range("R173").Value = 1.975
$ws.Range("S173").Value = 1.875
$ws.Range("T173").Value = 2.5
$ws.Range("U173").Value = 1.975
$ws.Range("V173").Value = 1.875
$ws.Range("W173").Value = 0.3999999999999999
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = -0.5
$ws.Range("AA173").Value = 0.4375
$ws.Range("AB173").Value = 0.9750000000000001
$ws.Range("AC173").Value = -1
